$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 11 (Bebidas) ---
$ws.Range("C11").Formula = "=B11*`$B4+B11"
$ws.Range("D11").Formula = "=C11*`$B4+C11"
$ws.Range("E11").Formula = "=D11*`$B4+D11"
$ws.Range("F11").Formula = "=SUM(B11:E11)"

# --- Row 12 (Condimentos) ---
$ws.Range("C12").Formula = "=B12*`$B5+B12"
$ws.Range("D12").Formula = "=C12*`$B5+C12"
$ws.Range("E12").Formula = "=D12*`$B5+D12"
$ws.Range("F12").Formula = "=SUM(B12:E12)"

# --- Row 13 (Lacteos) ---
$ws.Range("C13").Formula = "=B13*`$B6+B13"
$ws.Range("D13").Formula = "=C13*`$B6+C13"
$ws.Range("E13").Formula = "=D13*`$B6+D13"
$ws.Range("F13").Formula = "=SUM(B13:E13)"

# --- Row 14 (Frutas y Verduras) ---
$ws.Range("C14").Formula = "=B14*`$B7+B14"
$ws.Range("D14").Formula = "=C14*`$B7+C14"
$ws.Range("E14").Formula = "=D14*`$B7+D14"
$ws.Range("F14").Formula = "=SUM(B14:E14)"

# --- Row 15 totals (E15, F15 need new formula picking up extended range) ---
$ws.Range("E15").Formula = "=SUM(E11:E14)"
$ws.Range("F15").Formula = "=B15+C15+D15+E15"

# --- Row 16 (Gastos fijos) ---
$ws.Range("C16").Formula = "=B16*`$F4+B16"
$ws.Range("D16").Formula = "=C16*`$F4+C16"
$ws.Range("E16").Formula = "=D16*`$F4+D16"
$ws.Range("F16").Formula = "=SUM(B16:E16)"

# --- Row 17 (Gastos Variables) ---
$ws.Range("C17").Formula = "=B17*`$F5+B17"
$ws.Range("D17").Formula = "=C17*`$F5+C17"
$ws.Range("E17").Formula = "=D17*`$F5+D17"
$ws.Range("F17").Formula = "=SUM(B17:E17)"

# --- Row 18 totals (C18:F18) ---
$ws.Range("C18").Formula = "=SUM(C16:C17)"
$ws.Range("D18").Formula = "=SUM(D16:D17)"
$ws.Range("E18").Formula = "=SUM(E16:E17)"
$ws.Range("F18").Formula = "=SUM(F16:F17)"

# --- Row 19 totals (C19:F19) ---
$ws.Range("C19").Formula = "=C15-C18"
$ws.Range("D19").Formula = "=D15-D18"
$ws.Range("E19").Formula = "=E15-E18"
$ws.Range("F19").Formula = "=F15-F18"

# --- Row 20 cumulative totals (C20:F20) ---
$ws.Range("C20").Formula = "=B20+C19"
$ws.Range("D20").Formula = "=C20+D19"
$ws.Range("E20").Formula = "=D20+E19"
$ws.Range("F20").Formula = "=E20+F19"

# Update selection to match target (F19)
$ws.Range("F19").Select()
